# Update the "Metadata" worksheet (sheet1) to reflect the new ValueSet
# metadata (version bump, status change, new date, updated contact info,
# new Jurisdiction row, and the Immutable row moving down to make room).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$wsRx = $wb.Worksheets.Item("Include from RxNorm")

# --- Simple value edits (rows 3, 6, 8, 10 do not move) -------------------
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Make room for the new "Jurisdiction" row -----------------------------
# Row 16 gets the content that currently lives in row 15 ("Immutable").
# We copy format+value directly (bottom-most target first) so nothing is
# clobbered before it is copied onward. The destination is cleared first
# because copying an empty source cell on top of a non-empty destination
# would otherwise leave the destination's old value behind.
$ws.Range("A16:B16").ClearContents()
$ws.Range("A15:B15").Copy($ws.Range("A16:B16"))

# Row 15 becomes the current row 14 ("Copyright").
$ws.Range("A15:B15").ClearContents()
$ws.Range("A14:B14").Copy($ws.Range("A15:B15"))

# Row 14 becomes the current row 13 ("Purpose").
$ws.Range("A14:B14").ClearContents()
$ws.Range("A13:B13").Copy($ws.Range("A14:B14"))

# Row 13 becomes the current row 12 ("Description").
$ws.Range("A13:B13").ClearContents()
$ws.Range("A12:B12").Copy($ws.Range("A13:B13"))

# --- Row 11: still "Contact", but the value is now the named contact -----
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Row 12: brand-new "Jurisdiction" row with an empty value ------------
# Copy formatting from row 11 first so the new row matches the other data
# rows (border/fill/alignment), then set the label and leave the value
# blank.
$ws.Range("A12:B12").ClearContents()
$ws.Range("A11:B11").Copy($ws.Range("A12:B12"))
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# Row 16 keeps the "Immutable" / "BooleanType[null]" content/format that was
# copied down from the old row 15 above, so nothing further is required.
